$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J10").Value = 6
$ws.Range("J11").Value = 0
$ws.Range("J14").Value = "Herbst"
$ws.Range("J16").Value = "83-85 MPH"
$ws.Range("J17").Value = "SL,CB,FB,CH"

$ws.Range("J19").Value = 3
$ws.Range("M19").Value = "89.23 MPH"
$ws.Range("M21").Value = "22.91°"
$ws.Range("J23").Value = "Roblez"
$ws.Range("M23").Value = "Line Drive"
$ws.Range("M24").Value = "Double"
$ws.Range("J25").Value = "88-90 MPH"
$ws.Range("J26").Value = "CB,FB,CH"

$ws.Range("J28").Value = 2
$ws.Range("J29").Value = 0
$ws.Range("J32").Value = "Roblez"
$ws.Range("J33").Value = "Right"
$ws.Range("J34").Value = "88-90 MPH"
$ws.Range("J35").Value = "CB,FB,CH"

$ws.Range("J37").Value = 7
$ws.Range("J38").Value = 1
$ws.Range("J41").Value = "Plum"
$ws.Range("J43").Value = "84-86 MPH"
$ws.Range("J44").Value = "SL,FB,CH"

$ws.Range("J46").Value = 4
$ws.Range("M46").Value = "88.36 MPH"
$ws.Range("M48").Value = "49.18°"
$ws.Range("J50").Value = "Herbst"
$ws.Range("M50").Value = "Fly Ball"
$ws.Range("M51").Value = "Out"
$ws.Range("J52").Value = "83-85 MPH"
$ws.Range("J53").Value = "SL,CB,FB,CH"

$ws.Range("J61").Value = 8
$ws.Range("J62").Value = 2
$ws.Range("J65").Value = "Thompson"
$ws.Range("J66").Value = "Left"
$ws.Range("J67").Value = "84-84 MPH"
$ws.Range("J68").Value = "SL,FB,CH"
